# This workbook's "Sheet1" table (zz_Poker_PlayerInfo) is backed by a
# Power Query / SQL Server data source. The author's edit reflects a
# refreshed query result (new session stats for Jan instead of Dec),
# which changed several numeric cells in rows 2-10. Apply those
# refreshed values directly, then restore the sheet's last-known
# selection/active cell as recorded after the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 133
$ws.Range("E2").Value = 15
$ws.Range("J2").Value = 3.62

# Row 3
$ws.Range("D3").Value = 154
$ws.Range("E3").Value = 17
$ws.Range("I3").Value = 559
$ws.Range("J3").Value = 3.63

# Row 4
$ws.Range("D4").Value = 203
$ws.Range("E4").Value = 21
$ws.Range("I4").Value = 748
$ws.Range("J4").Value = 3.68

# Row 5
$ws.Range("D5").Value = 115
$ws.Range("E5").Value = 14
$ws.Range("I5").Value = 419
$ws.Range("J5").Value = 3.64

# Row 7
$ws.Range("D7").Value = 155
$ws.Range("E7").Value = 21
$ws.Range("I7").Value = 669
$ws.Range("J7").Value = 4.32

# Row 8
$ws.Range("D8").Value = 208
$ws.Range("E8").Value = 21
$ws.Range("I8").Value = 720

# Row 10
$ws.Range("D10").Value = 211
$ws.Range("E10").Value = 21
$ws.Range("I10").Value = 880
$ws.Range("J10").Value = 4.17

# Update the saved selection/active cell on the sheet.
$ws.Range("I35").Select()
